$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.719.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.600.77"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.14"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.246"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.66"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.639.85"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.13"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.695.15"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "210.33"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.29"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.95"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.99"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.09"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.35"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.25"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.97"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.292.46"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.47"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.602"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.17"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +14.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0169"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.20"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.11"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.738.63"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.65"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.02%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.78%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.38"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.77%  "
